$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'69.077.66"
$ws.Range("E2").Value = "  +1.87%  "
$ws.Range("D3").Value = "'3.768.48"
$ws.Range("E3").Value = "  -0.45%  "
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'625.55"
$ws.Range("E5").Value = "  +3.87%  "
$ws.Range("D6").Value = "'165.53"
$ws.Range("E6").Value = "  +1.55%  "
$ws.Range("D7").Value = "'3.767.66"
$ws.Range("E7").Value = "  -0.40%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E9").Value = "  +1.50%  "
$ws.Range("E10").Value = "  +1.45%  "
$ws.Range("D11").Value = "'0.459"
$ws.Range("E11").Value = "  +3.15%  "
$ws.Range("D12").Value = "'6.74"
$ws.Range("E12").Value = "  -0.04%  "
$ws.Range("D13").Value = "'0.0000246"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "'35.65"
$ws.Range("E14").Value = "  +1.63%  "
$ws.Range("D15").Value = "'4.400.45"
$ws.Range("E15").Value = "  -0.46%  "
$ws.Range("D16").Value = "'3.643.25"
$ws.Range("E16").Value = "  -3.37%  "
$ws.Range("D17").Value = "'69.074.82"
$ws.Range("E17").Value = "  +1.88%  "
$ws.Range("D18").Value = "'17.66"
$ws.Range("E18").Value = "  -2.73%  "
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "'7.06"
$ws.Range("E20").Value = "  +0.80%  "
$ws.Range("D21").Value = "'467.12"
$ws.Range("E21").Value = "  +2.20%  "
$ws.Range("D22").Value = "'9.56"
$ws.Range("E22").Value = "  +1.33%  "
$ws.Range("D23").Value = "'0.706"
$ws.Range("E23").Value = "  +2.49%  "
$ws.Range("E24").Value = "  +3.28%  "
$ws.Range("D25").Value = "'83.24"
$ws.Range("E25").Value = "  +0.31%  "
$ws.Range("D26").Value = "'12.02"
$ws.Range("E26").Value = "  +1.62%  "
$ws.Range("E27").Value = "  +3.80%  "
$ws.Range("D28").Value = "'10.03"
$ws.Range("E28").Value = "  +1.41%  "
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").Value = "'3.917.86"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("E31").Value = "  +2.74%  "
$ws.Range("E32").Value = "  +2.57%  "
$ws.Range("E33").Value = "  -0.56%  "
$ws.Range("D34").Value = "'28.74"
$ws.Range("E34").Value = "  -0.44%  "
$ws.Range("D35").Value = "'0.173"
$ws.Range("E35").Value = "  +20.11%  "
$ws.Range("E36").Value = "  +0.16%  "
$ws.Range("D37").Value = "'3.719.84"
$ws.Range("E37").Value = "  -0.40%  "
$ws.Range("D38").Value = "'8.94"
$ws.Range("E38").Value = "  +0.42%  "
$ws.Range("E39").Value = "  +2.26%  "
$ws.Range("E41").Value = "  +0.60%  "
$ws.Range("D42").Value = "'0.967"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "'153.54"
$ws.Range("E45").Value = "  +1.19%  "
$ws.Range("D46").Value = "'43.16"
$ws.Range("E46").Value = "  -1.31%  "
$ws.Range("D47").Value = "'0.296"
$ws.Range("E47").Value = "  +0.75%  "
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "'1.90"
$ws.Range("E49").Value = "  +4.12%  "
$ws.Range("D50").Value = "'8.39"
$ws.Range("E50").Value = "  +1.59%  "
$ws.Range("E51").Value = "  +0.16%  "
